$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "browserId"
$ws.Range("B1").Value = "parentPage"
$ws.Range("C1").Value = "fullMediaUrl"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "src"
$ws.Range("F1").Value = "status"
$ws.Range("G1").Value = "ok"
$ws.Range("H1").Value = "error"

# Row 2
$ws.Range("A2").Value = "B2-W1"
$ws.Range("B2").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C2").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Acquired%20Brain%20Injury%20Documentation%20Form%20SUNY%20Geneseo.pdf"
$ws.Range("D2").Value = "pdf"
$ws.Range("E2").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Acquired%20Brain%20Injury%20Documentation%20Form%20SUNY%20Geneseo.pdf"
$ws.Range("F2").Value = 404
$ws.Range("G2").Value = $false

# Row 3
$ws.Range("A3").Value = "B2-W1"
$ws.Range("B3").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C3").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Attention%20Deficit%20Documentation.pdf"
$ws.Range("D3").Value = "pdf"
$ws.Range("E3").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Attention%20Deficit%20Documentation.pdf"
$ws.Range("F3").Value = 404
$ws.Range("G3").Value = $false

# Row 4
$ws.Range("A4").Value = "B2-W1"
$ws.Range("B4").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C4").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Autism%20Spectrum%20Documentation%20Form.pdf"
$ws.Range("D4").Value = "pdf"
$ws.Range("E4").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Autism%20Spectrum%20Documentation%20Form.pdf"
$ws.Range("F4").Value = 404
$ws.Range("G4").Value = $false

# Row 5
$ws.Range("A5").Value = "B2-W1"
$ws.Range("B5").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C5").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Deaf%20and%20HH%20documentation.pdf"
$ws.Range("D5").Value = "pdf"
$ws.Range("E5").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Deaf%20and%20HH%20documentation.pdf"
$ws.Range("F5").Value = 404
$ws.Range("G5").Value = $false

# Row 6
$ws.Range("A6").Value = "B2-W1"
$ws.Range("B6").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C6").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Learning%20Disabilities%20Documentation%20Guidelines.pdf"
$ws.Range("D6").Value = "pdf"
$ws.Range("E6").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Learning%20Disabilities%20Documentation%20Guidelines.pdf"
$ws.Range("F6").Value = 404
$ws.Range("G6").Value = $false

# Row 7
$ws.Range("A7").Value = "B2-W1"
$ws.Range("B7").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C7").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Chronic%20Health%20Conditions%20Documentation.pdf"
$ws.Range("D7").Value = "pdf"
$ws.Range("E7").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Chronic%20Health%20Conditions%20Documentation.pdf"
$ws.Range("F7").Value = 404
$ws.Range("G7").Value = $false

# Row 8
$ws.Range("A8").Value = "B2-W1"
$ws.Range("B8").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C8").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/CONCUSSION%20EVALUATION.pdf"
$ws.Range("D8").Value = "pdf"
$ws.Range("E8").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/CONCUSSION%20EVALUATION.pdf"
$ws.Range("F8").Value = 404
$ws.Range("G8").Value = $false

# Row 9
$ws.Range("A9").Value = "B2-W1"
$ws.Range("B9").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C9").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Mobility%20Impairment%20Documentation.pdf"
$ws.Range("D9").Value = "pdf"
$ws.Range("E9").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Mobility%20Impairment%20Documentation.pdf"
$ws.Range("F9").Value = 404
$ws.Range("G9").Value = $false

# Row 10
$ws.Range("A10").Value = "B2-W1"
$ws.Range("B10").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C10").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Neurological%20Disorders%20Documentation.pdf"
$ws.Range("D10").Value = "pdf"
$ws.Range("E10").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Neurological%20Disorders%20Documentation.pdf"
$ws.Range("F10").Value = 404
$ws.Range("G10").Value = $false

# Row 11
$ws.Range("A11").Value = "B2-W1"
$ws.Range("B11").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C11").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Psychological%20Conditions%20Documentation.pdf"
$ws.Range("D11").Value = "pdf"
$ws.Range("E11").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Psychological%20Conditions%20Documentation.pdf"
$ws.Range("F11").Value = 404
$ws.Range("G11").Value = $false

# Row 12
$ws.Range("A12").Value = "B2-W1"
$ws.Range("B12").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C12").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Temporary%20Impairments%20Documentation.pdf"
$ws.Range("D12").Value = "pdf"
$ws.Range("E12").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Temporary%20Impairments%20Documentation.pdf"
$ws.Range("F12").Value = 404
$ws.Range("G12").Value = $false

# Row 13
$ws.Range("A13").Value = "B2-W1"
$ws.Range("B13").Value = "http://localhost/sunny/accessibility-office/documentation-criteria"
$ws.Range("C13").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Visual%20Impairment%20Documentation.pdf"
$ws.Range("D13").Value = "pdf"
$ws.Range("E13").Value = "http://localhost/sunny/sites/default/files/sites/dean_office/Visual%20Impairment%20Documentation.pdf"
$ws.Range("F13").Value = 404
$ws.Range("G13").Value = $false

# Row 14
$ws.Range("A14").Value = "B1-W2"
$ws.Range("B14").Value = "http://localhost/sunny/aac/forms-and-documents"
$ws.Range("C14").Value = "http://localhost/sites/default/files/users/1120/Continuing%20Recognition%20Form.pdf"
$ws.Range("D14").Value = "pdf"
$ws.Range("E14").Value = "http://localhost/sites/default/files/users/1120/Continuing%20Recognition%20Form.pdf"
$ws.Range("F14").Value = 404
$ws.Range("G14").Value = $false

